# Updates cryptos list values (price & 1h volume) per the scraped diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.913.26"
$ws.Range("D3").Value = "3.116.42"
$ws.Range("E3").Value = "  +2.99%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'585.49"
$ws.Range("E5").Value = "  +3.39%  "
$ws.Range("D6").Value = "'144.63"
$ws.Range("E6").Value = "  +2.43%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.109.00"
$ws.Range("E8").Value = "  +3.04%  "
$ws.Range("D9").Value = "'0.530"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("D10").Value = "'0.150"
$ws.Range("E10").Value = "  +11.24%  "
$ws.Range("D11").Value = "'5.70"
$ws.Range("E11").Value = "  +7.42%  "
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "  +5.72%  "
$ws.Range("D14").Value = "'35.50"
$ws.Range("E14").Value = "  +3.64%  "
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "3.629.98"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "3.109.68"
$ws.Range("E18").Value = "  +2.85%  "
$ws.Range("D19").Value = "62.838.04"
$ws.Range("E19").Value = "  +4.73%  "
$ws.Range("D20").Value = "'465.49"
$ws.Range("E20").Value = "  +6.03%  "
$ws.Range("D21").Value = "'14.09"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("D22").Value = "'0.729"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "'7.54"
$ws.Range("E23").Value = "  +5.26%  "
$ws.Range("D24").Value = "'13.37"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").Value = "'82.04"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").Value = "'2.68"
$ws.Range("E28").Value = "  +4.53%  "
$ws.Range("D29").Value = "'8.31"
$ws.Range("E29").Value = "  +5.29%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'6.84"
$ws.Range("E31").Value = "  +7.96%  "
$ws.Range("D32").Value = "'26.99"
$ws.Range("E32").Value = "  +3.19%  "
$ws.Range("D33").Value = "'0.111"
$ws.Range("E33").Value = "  +8.73%  "
$ws.Range("D34").Value = "0.0₃0842"
$ws.Range("E34").Value = "  +6.39%  "
$ws.Range("D35").Value = "'2.37"
$ws.Range("E35").Value = "  +10.81%  "
$ws.Range("E36").Value = "  +3.52%  "
$ws.Range("D37").Value = "'6.05"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").Value = "'3.19"
$ws.Range("E38").Value = "  +13.99%  "
$ws.Range("D39").Value = "'51.11"
$ws.Range("E39").Value = "  +3.84%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'433.03"
$ws.Range("E40").Value = "  +6.44%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'8.80"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D42").Value = "2.931.23"
$ws.Range("E42").Value = "  +5.08%  "
$ws.Range("D43").Value = "'0.0369"
$ws.Range("E43").Value = "  +3.54%  "
$ws.Range("E44").Value = "  +8.98%  "
$ws.Range("D45").Value = "'0.111"
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("E46").Value = "  +6.23%  "
$ws.Range("D47").Value = "'35.22"
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("D49").Value = "'123.37"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").Value = "'24.73"
$ws.Range("E51").Value = "  +4.31%  "
